$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching style/format of existing header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill H2:H14 with 0 for each data row
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
}
